$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2078.111
$ws.Range("I11").Value = 2078.111
$ws.Range("K11").Value = 2078.111
$ws.Range("M11").Value = -1938.111
$ws.Range("H28").Value = 1123
$ws.Range("I28").Value = 635.3
$ws.Range("K28").Value = 635.3
$ws.Range("M28").Value = -150.3
$ws.Range("H29").Value = 777
$ws.Range("I29").Value = 777
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2331
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2050
$ws.Range("N29").Value = ""
$ws.Range("H41").Value = 542.0769
$ws.Range("I41").Value = 630.8889
$ws.Range("J41").Value = 342.25
$ws.Range("K41").Value = 630.8889
$ws.Range("L41").Value = 342.25
$ws.Range("M41").Value = -190.8889
$ws.Range("N41").Value = -1222.25
$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 6000
$ws.Range("K62").Value = 6000
$ws.Range("M62").Value = -5376
$ws.Range("H64").Value = 4801
$ws.Range("I64").Value = 3666.3333
$ws.Range("K64").Value = 3666.3333
$ws.Range("M64").Value = -3418.3333
$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 6000
$ws.Range("K65").Value = 30000
$ws.Range("M65").Value = -26880
$ws.Range("H67").Value = 4801
$ws.Range("I67").Value = 3666.3333
$ws.Range("K67").Value = 3666.3333
$ws.Range("M67").Value = -2808.3333
$ws.Range("H80").Value = 4999.5
$ws.Range("I80").Value = 4999
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 14997
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = -13999
$ws.Range("N80").Value = -16996
$ws.Range("H83").Value = 4999.5
$ws.Range("I83").Value = 4999
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 44991
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -39999
$ws.Range("N83").Value = -54984
$ws.Range("H98").Value = 464.33334
$ws.Range("I98").Value = 464.33334
$ws.Range("K98").Value = 464.33334
$ws.Range("M98").Value = 1033.66666
$ws.Range("H106").Value = 12741.571
$ws.Range("I106").Value = 9865.166999999999
$ws.Range("K106").Value = 9865.166999999999
$ws.Range("M106").Value = -9234.166999999999
$ws.Range("H122").Value = 464.33334
$ws.Range("I122").Value = 464.33334
$ws.Range("K122").Value = 1393.00002
$ws.Range("M122").Value = 1056.99998
$ws.Range("H130").Value = 27744.5
$ws.Range("I130").Value = 29709
$ws.Range("J130").Value = 25780
$ws.Range("K130").Value = 29709
$ws.Range("L130").Value = 25780
$ws.Range("M130").Value = -24689
$ws.Range("N130").Value = -35820

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("H32").Value = 1881.6207
$ws.Range("I32").Value = 1835.8148
$ws.Range("K32").Value = 1835.8148
$ws.Range("M32").Value = -1548.8148
$ws.Range("H45").Value = 2759.4
$ws.Range("I45").Value = 2466.6667
$ws.Range("K45").Value = 2466.6667
$ws.Range("M45").Value = -2089.6667
$ws.Range("H63").Value = 3731.1794
$ws.Range("I63").Value = 2412.5757
$ws.Range("J63").Value = 10983.5
$ws.Range("K63").Value = 2412.5757
$ws.Range("L63").Value = 10983.5
$ws.Range("M63").Value = -1726.5757
$ws.Range("N63").Value = -12355.5
$ws.Range("H66").Value = 3731.1794
$ws.Range("I66").Value = 2412.5757
$ws.Range("J66").Value = 10983.5
$ws.Range("K66").Value = 12062.8785
$ws.Range("L66").Value = 54917.5
$ws.Range("M66").Value = -8630.878499999999
$ws.Range("N66").Value = -61781.5
$ws.Range("H97").Value = 500
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 500
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -4
$ws.Range("N97").Value = ""
$ws.Range("H132").Value = 3481.4546
$ws.Range("I132").Value = 1899.5
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 5698.5
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -3168.5
$ws.Range("N132").Value = -16559

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7487.25
$ws.Range("J80").Value = 14775.75
$ws.Range("L80").Value = 14775.75
$ws.Range("N80").Value = -16771.75
$ws.Range("H83").Value = 7487.25
$ws.Range("J83").Value = 14775.75
$ws.Range("L83").Value = 73878.75
$ws.Range("N83").Value = -83862.75
$ws.Range("H107").Value = 3477.2856
$ws.Range("I107").Value = 3477.2856
$ws.Range("K107").Value = 3477.2856
$ws.Range("M107").Value = -1557.2856
$ws.Range("H110").Value = 60001
$ws.Range("J110").Value = 60001
$ws.Range("L110").Value = 60001
$ws.Range("N110").Value = -68181

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4797
$ws.Range("I31").Value = 4797
$ws.Range("K31").Value = 4797
$ws.Range("M31").Value = -4502
$ws.Range("H34").Value = 4797
$ws.Range("I34").Value = 4797
$ws.Range("K34").Value = 4797
$ws.Range("M34").Value = -4595
$ws.Range("H134").Value = 4293.6665
$ws.Range("I134").Value = 3950
$ws.Range("J134").Value = 6012
$ws.Range("K134").Value = 11850
$ws.Range("L134").Value = 18036
$ws.Range("M134").Value = -9315
$ws.Range("N134").Value = -23106

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 397.5
$ws.Range("I23").Value = 299
$ws.Range("J23").Value = 693
$ws.Range("K23").Value = 897
$ws.Range("L23").Value = 2079
$ws.Range("M23").Value = -662
$ws.Range("N23").Value = -2549
$ws.Range("H50").Value = 1843.5714
$ws.Range("I50").Value = 192.4
$ws.Range("K50").Value = 577.2
$ws.Range("M50").Value = -96.20000000000005
$ws.Range("H53").Value = 1843.5714
$ws.Range("I53").Value = 192.4
$ws.Range("K53").Value = 577.2
$ws.Range("M53").Value = -96.20000000000005
$ws.Range("H131").Value = 1986.0526
$ws.Range("I131").Value = 1595.0769
$ws.Range("K131").Value = 4785.2307
$ws.Range("M131").Value = 254.7692999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6061.909
$ws.Range("I70").Value = 4459.75
$ws.Range("K70").Value = 4459.75
$ws.Range("M70").Value = -4189.75
$ws.Range("H73").Value = 6061.909
$ws.Range("I73").Value = 4459.75
$ws.Range("K73").Value = 4459.75
$ws.Range("M73").Value = -3523.75
$ws.Range("H102").Value = 4675
$ws.Range("I102").Value = 5310
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 5310
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -3688
$ws.Range("N102").Value = -4744
$ws.Range("H126").Value = 2731.5454
$ws.Range("I126").Value = 2803.3
$ws.Range("K126").Value = 8409.900000000001
$ws.Range("M126").Value = -5939.900000000001
$ws.Range("H132").Value = 3923.4614
$ws.Range("I132").Value = 3725.05
$ws.Range("J132").Value = 4584.8335
$ws.Range("K132").Value = 11175.15
$ws.Range("L132").Value = 13754.5005
$ws.Range("M132").Value = -8645.150000000001
$ws.Range("N132").Value = -18814.5005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H22").Value = 1598.3636
$ws.Range("I22").Value = 1172.6364
$ws.Range("J22").Value = 2024.091
$ws.Range("K22").Value = 1172.6364
$ws.Range("L22").Value = 2024.091
$ws.Range("M22").Value = -877.6364000000001
$ws.Range("N22").Value = -2614.091
$ws.Range("H27").Value = 1598.3636
$ws.Range("I27").Value = 1172.6364
$ws.Range("J27").Value = 2024.091
$ws.Range("K27").Value = 1172.6364
$ws.Range("L27").Value = 2024.091
$ws.Range("M27").Value = -1065.6364
$ws.Range("N27").Value = -2238.091
$ws.Range("H40").Value = 4599.4443
$ws.Range("I40").Value = 4599.4443
$ws.Range("K40").Value = 4599.4443
$ws.Range("M40").Value = -4463.4443
$ws.Range("H55").Value = 1176.1428
$ws.Range("I55").Value = 945.8
$ws.Range("J55").Value = 1304.1111
$ws.Range("K55").Value = 945.8
$ws.Range("L55").Value = 1304.1111
$ws.Range("M55").Value = -772.8
$ws.Range("N55").Value = -1650.1111
$ws.Range("H68").Value = 6736
$ws.Range("J68").Value = 6450
$ws.Range("L68").Value = 6450
$ws.Range("N68").Value = -7948
$ws.Range("H71").Value = 6736
$ws.Range("J71").Value = 6450
$ws.Range("L71").Value = 32250
$ws.Range("N71").Value = -39738
$ws.Range("H122").Value = 6114.7144
$ws.Range("I122").Value = 6717.1665
$ws.Range("K122").Value = 20151.4995
$ws.Range("M122").Value = -17701.4995
$ws.Range("H130").Value = 84950
$ws.Range("J130").Value = 84950
$ws.Range("L130").Value = 84950
$ws.Range("N130").Value = -94990
$ws.Range("H132").Value = 2188.8696
$ws.Range("I132").Value = 2150
$ws.Range("J132").Value = 2448
$ws.Range("K132").Value = 6450
$ws.Range("L132").Value = 7344
$ws.Range("M132").Value = -3920
$ws.Range("N132").Value = -12404

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 19999.5
$ws.Range("J82").Value = 19999.5
$ws.Range("L82").Value = 19999.5
$ws.Range("N82").Value = -20765.5
$ws.Range("H85").Value = 19999.5
$ws.Range("J85").Value = 19999.5
$ws.Range("L85").Value = 19999.5
$ws.Range("N85").Value = -22651.5
$ws.Range("H113").Value = 833.8461
$ws.Range("I113").Value = 828.3333
$ws.Range("K113").Value = 2484.9999
$ws.Range("M113").Value = -314.9998999999998
$ws.Range("H122").Value = 8002
$ws.Range("I122").Value = 9336
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 28008
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -25558
$ws.Range("N122").Value = -16900
